$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 742.2857
$ws.Range("J2").Value = 736.75
$ws.Range("L2").Value = 736.75
$ws.Range("N2").Value = -962.75
$ws.Range("H38").Value = 3032.7058
$ws.Range("I38").Value = 49.727272
$ws.Range("J38").Value = 8501.5
$ws.Range("K38").Value = 149.181816
$ws.Range("L38").Value = 25504.5
$ws.Range("M38").Value = 222.818184
$ws.Range("N38").Value = -26248.5
$ws.Range("H58").Value = 1283.5
$ws.Range("I58").Value = 372.5
$ws.Range("J58").Value = 4016.5
$ws.Range("K58").Value = 1117.5
$ws.Range("L58").Value = 12049.5
$ws.Range("M58").Value = -967.5
$ws.Range("N58").Value = -12349.5
$ws.Range("H76").Value = 7438.5864
$ws.Range("I76").Value = 7240.5
$ws.Range("K76").Value = 7240.5
$ws.Range("M76").Value = -6925.5
$ws.Range("H79").Value = 7438.5864
$ws.Range("I79").Value = 7240.5
$ws.Range("K79").Value = 7240.5
$ws.Range("M79").Value = -6148.5
$ws.Range("H86").Value = 4373.6313
$ws.Range("I86").Value = 4905.5
$ws.Range("J86").Value = 3782.6667
$ws.Range("K86").Value = 4905.5
$ws.Range("L86").Value = 3782.6667
$ws.Range("M86").Value = -3782.5
$ws.Range("N86").Value = -6028.6667
$ws.Range("H89").Value = 4373.6313
$ws.Range("I89").Value = 4905.5
$ws.Range("J89").Value = 3782.6667
$ws.Range("K89").Value = 24527.5
$ws.Range("L89").Value = 18913.3335
$ws.Range("M89").Value = -18911.5
$ws.Range("N89").Value = -30145.3335
$ws.Range("H96").Value = 1192.2727
$ws.Range("I96").Value = 764.3333
$ws.Range("J96").Value = 1705.8
$ws.Range("K96").Value = 2292.9999
$ws.Range("L96").Value = 5117.4
$ws.Range("M96").Value = -919.9998999999998
$ws.Range("N96").Value = -7863.4
$ws.Range("H107").Value = 1804.9375
$ws.Range("I107").Value = 1634.3572
$ws.Range("J107").Value = 2999
$ws.Range("K107").Value = 1634.3572
$ws.Range("L107").Value = 2999
$ws.Range("M107").Value = 285.6428000000001
$ws.Range("N107").Value = -6839

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15774.353
$ws.Range("J45").Value = 5687.5
$ws.Range("L45").Value = 5687.5
$ws.Range("N45").Value = -6441.5
$ws.Range("H63").Value = 8319.333000000001
$ws.Range("I63").Value = 7148.1113
$ws.Range("K63").Value = 7148.1113
$ws.Range("M63").Value = -6462.1113
$ws.Range("H66").Value = 8319.333000000001
$ws.Range("I66").Value = 7148.1113
$ws.Range("K66").Value = 35740.5565
$ws.Range("M66").Value = -32308.5565
$ws.Range("H88").Value = 2548.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2548.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2548.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3360.5
$ws.Range("H91").Value = 2548.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2548.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2548.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5356.5
$ws.Range("H122").Value = 2045.9584
$ws.Range("I122").Value = 1890.4762
$ws.Range("K122").Value = 5671.4286
$ws.Range("M122").Value = -3221.4286
$ws.Range("H133").Value = 43657.91
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4156.9165
$ws.Range("I86").Value = 3116.3333
$ws.Range("J86").Value = 5197.5
$ws.Range("K86").Value = 3116.3333
$ws.Range("L86").Value = 5197.5
$ws.Range("M86").Value = -1993.3333
$ws.Range("N86").Value = -7443.5
$ws.Range("H89").Value = 4156.9165
$ws.Range("I89").Value = 3116.3333
$ws.Range("J89").Value = 5197.5
$ws.Range("K89").Value = 15581.6665
$ws.Range("L89").Value = 25987.5
$ws.Range("M89").Value = -9965.666499999999
$ws.Range("N89").Value = -37219.5
$ws.Range("H134").Value = 2655.7778
$ws.Range("I134").Value = 2655.7778
$ws.Range("K134").Value = 7967.3334
$ws.Range("M134").Value = -5432.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9049.947
$ws.Range("I31").Value = 11768.308
$ws.Range("K31").Value = 11768.308
$ws.Range("M31").Value = -11473.308
$ws.Range("H34").Value = 9049.947
$ws.Range("I34").Value = 11768.308
$ws.Range("K34").Value = 11768.308
$ws.Range("M34").Value = -11566.308
$ws.Range("H105").Value = 1131.7273
$ws.Range("I105").Value = 1156.25
$ws.Range("K105").Value = 1156.25
$ws.Range("M105").Value = 590.75
$ws.Range("H107").Value = 1139.1482
$ws.Range("I107").Value = 1125.3334
$ws.Range("K107").Value = 1125.3334
$ws.Range("M107").Value = 794.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 10250.2
$ws.Range("I57").Value = 3500.75
$ws.Range("K57").Value = 10502.25
$ws.Range("M57").Value = -9943.25
$ws.Range("H80").Value = 4000
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 9000
$ws.Range("N80").Value = -10872
$ws.Range("H83").Value = 4000
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 27000
$ws.Range("N83").Value = -36360
$ws.Range("H121").Value = 633.4375
$ws.Range("J121").Value = 1896.6666
$ws.Range("L121").Value = 5689.9998
$ws.Range("N121").Value = -8309.9998
$ws.Range("H128").Value = 165999.5
$ws.Range("I128").Value = 165999.5
$ws.Range("K128").Value = 497998.5
$ws.Range("M128").Value = -493018.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7882.5
$ws.Range("I113").Value = 8745
$ws.Range("J113").Value = 1845
$ws.Range("K113").Value = 8745
$ws.Range("L113").Value = 1845
$ws.Range("M113").Value = -6575
$ws.Range("N113").Value = -6185

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3900
$ws.Range("I100").Value = 2062.5
$ws.Range("K100").Value = 2062.5
$ws.Range("M100").Value = -1521.5
$ws.Range("H138").Value = 106285
$ws.Range("J138").Value = 106285
$ws.Range("L138").Value = 106285
$ws.Range("N138").Value = -116565

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 30000
$ws.Range("I103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("K103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("M103").Value = -28828
$ws.Range("N103").Value = -32344
$ws.Range("H122").Value = 4752.8423
$ws.Range("I122").Value = 4135.5293
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 12406.5879
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -9956.5879
$ws.Range("N122").Value = -34900
$ws.Range("H136").Value = 1723.5
$ws.Range("I136").Value = 1139.7273
$ws.Range("K136").Value = 3419.1819
$ws.Range("M136").Value = -869.1819
